$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.670.79'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = '  +0.47%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.149.08'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = '  +1.56%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '570.62'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +1.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.52'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +4.16%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.148.05'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  +1.53%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.528'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  +3.68%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.162'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +4.91%  '
$ws.Range("E11").Value = '  +0.14%  '
$ws.Range("E12").Value = '  +7.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000259'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  +12.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '38.22'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  +8.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.665.26'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  +1.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.771.57'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  +0.60%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.24'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  +7.03%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.150.38'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  +1.75%  '
$ws.Range("E19").Value = '  +0.30%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '516.17'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  +6.56%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.98'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +6.83%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.739'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  +8.92%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.31'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +6.81%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.87'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  +3.78%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.33'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  +4.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("E27").Value = '  +4.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.76'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  +8.99%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.19'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  +6.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '27.95'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  +5.65%  '
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("B32").Value = 'Mantle'
$ws.Range("C32").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.19'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  +4.25%  '
$ws.Range("B33").Value = 'Stacks'
$ws.Range("C33").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.68'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  +8.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.14'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +8.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.60'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  +5.67%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.71'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  +0.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '488.18'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  +9.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0866'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  +5.91%  '
$ws.Range("E39").Value = '  +3.38%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.98'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  -0.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.113.51'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  +4.86%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.68'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +5.28%  '
$ws.Range("E43").Value = '  +5.74%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.294'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  +12.29%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.45'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  +14.53%  '
$ws.Range("E46").Value = '  +4.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₃0577'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  +11.54%  '
$ws.Range("E49").Value = '  +2.86%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.30'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +10.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '119.14'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  +0.52%  '
